{"js": "// Update the \"Genk\u00e4rom\u00e5l\" (counterclaim) section: raise the total from\n// 150 000 kr to 400 000 kr, bump two of the itemised amounts, reword the\n// \"economic damage\" bullet, add a new \"Sveda och v\u00e4rk\" bullet, and update\n// the cross-reference filename.\n\n// 1) The heading and the intro sentence both contain the literal phrase\n//    \"150 000 kr\" \u2014 replace every occurrence with \"400 000 kr\".\nconst totalMatches = context.document.body.search(\"150 000 kr\", { matchCase: true });\ntotalMatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < totalMatches.items.length; i++) {\n  totalMatches.items[i].insertText(\"400 000 kr\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Bump the \"Kr\u00e4nkningsers\u00e4ttning\" bullet amount: 50 000 kr -> 150 000 kr.\nconst krankningMatches = context.document.body.search(\n  \"50 000 kr \u2013 Kr\u00e4nkningsers\u00e4ttning f\u00f6r falsk v\u00e5ldt\u00e4ktsanm\u00e4lan (2 kap 3 \u00a7 SkL, BrB 15:7)\",\n  { matchCase: true }\n);\nkrankningMatches.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < krankningMatches.items.length; i++) {\n  krankningMatches.items[i].insertText(\n    \"150 000 kr \u2013 Kr\u00e4nkningsers\u00e4ttning f\u00f6r falsk v\u00e5ldt\u00e4ktsanm\u00e4lan (2 kap 3 \u00a7 SkL, BrB 15:7)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 3) Bump the \"Skadest\u00e5nd f\u00f6r r\u00e4tteg\u00e5ngsmissbruk\" bullet: 40 000 kr -> 100 000 kr.\nconst skadestandMatches = context.document.body.search(\n  \"40 000 kr \u2013 Skadest\u00e5nd f\u00f6r r\u00e4tteg\u00e5ngsmissbruk (2 kap 2 \u00a7 SkL, jfr BrB 15:2)\",\n  { matchCase: true }\n);\nskadestandMatches.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < skadestandMatches.items.length; i++) {\n  skadestandMatches.items[i].insertText(\n    \"100 000 kr \u2013 Skadest\u00e5nd f\u00f6r r\u00e4tteg\u00e5ngsmissbruk (2 kap 2 \u00a7 SkL, jfr BrB 15:2)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 4) Reword + bump the economic-damage bullet: 25 000 kr -> 75 000 kr, and\n//    extend the description with \"f\u00f6rlorad inkomst, resekostnader\".\nconst ekonomiskMatches = context.document.body.search(\n  \"25 000 kr \u2013 Ekonomisk skada (\u00f6vers\u00e4ttningskostnader, f\u00f6rlorad arbetstid m.m.)\",\n  { matchCase: true }\n);\nekonomiskMatches.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < ekonomiskMatches.items.length; i++) {\n  ekonomiskMatches.items[i].insertText(\n    \"75 000 kr \u2013 Ekonomisk skada (\u00f6vers\u00e4ttningskostnader, f\u00f6rlorad inkomst, resekostnader m.m.)\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n\n// 5) Insert a brand-new bullet for \"Sveda och v\u00e4rk\" right after the\n//    (now reworded) economic-damage bullet, before the closing sentence.\nconst newBulletAnchor = context.document.body.search(\n  \"75 000 kr \u2013 Ekonomisk skada (\u00f6vers\u00e4ttningskostnader, f\u00f6rlorad inkomst, resekostnader m.m.)\",\n  { matchCase: true }\n);\nnewBulletAnchor.load(\"text\");\nawait context.sync();\n\nif (newBulletAnchor.items.length > 0) {\n  const anchorParagraph = newBulletAnchor.items[0].paragraphs.getFirst();\n  const newParagraph = anchorParagraph.insertParagraph(\n    \"40 000 kr \u2013 Sveda och v\u00e4rk / psykiskt lidande (5 kap 1 \u00a7 SkL)\",\n    Word.InsertLocation.after\n  );\n  newParagraph.style = \"List Bullet\";\n}\nawait context.sync();\n\n// 6) Update the cross-reference filename in the closing sentence.\nconst fileRefMatches = context.document.body.search(\n  \"08_GENKAROMAL_150000kr.docx\",\n  { matchCase: true }\n);\nfileRefMatches.load(\"text\");\nawait context.sync();\nfor (let i = 0; i < fileRefMatches.items.length; i++) {\n  fileRefMatches.items[i].insertText(\n    \"08_GENKAROMAL_400000kr.docx\",\n    Word.InsertLocation.replace\n  );\n}\nawait context.sync();\n", "ps1": "# Update the \"Genk\u00e4rom\u00e5l\" (counterclaim) section: raise the total from\n# 150 000 kr to 400 000 kr, bump two of the itemised amounts, reword the\n# \"economic damage\" bullet, add a new \"Sveda och v\u00e4rk\" bullet, and update\n# the cross-reference filename.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $findText\n    $rng.Find.Replacement.Text = $replaceText\n    $rng.Find.Execute(\n        [ref]$findText,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        [ref]$replaceText,\n        2\n    )\n}\n\nfunction Find-ParagraphIndex($searchText) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $searchText\n    $found = $rng.Find.Execute()\n    if (-not $found) {\n        return -1\n    }\n    $paras = $d.Paragraphs\n    for ($i = 1; $i -le $paras.Count; $i++) {\n        $p = $paras.Item($i)\n        if ($p.Range.Start -le $rng.Start -and $p.Range.End -ge $rng.End) {\n            return $i\n        }\n    }\n    return -1\n}\n\n# 1) Heading + intro sentence: \"150 000 kr\" -> \"400 000 kr\" (both occurrences).\nReplace-All \"150 000 kr\" \"400 000 kr\"\n\n# 2) \"Kr\u00e4nkningsers\u00e4ttning\" bullet: 50 000 kr -> 150 000 kr.\nReplace-All `\n    \"50 000 kr \u2013 Kr\u00e4nkningsers\u00e4ttning f\u00f6r falsk v\u00e5ldt\u00e4ktsanm\u00e4lan (2 kap 3 \u00a7 SkL, BrB 15:7)\" `\n    \"150 000 kr \u2013 Kr\u00e4nkningsers\u00e4ttning f\u00f6r falsk v\u00e5ldt\u00e4ktsanm\u00e4lan (2 kap 3 \u00a7 SkL, BrB 15:7)\"\n\n# 3) \"Skadest\u00e5nd f\u00f6r r\u00e4tteg\u00e5ngsmissbruk\" bullet: 40 000 kr -> 100 000 kr.\nReplace-All `\n    \"40 000 kr \u2013 Skadest\u00e5nd f\u00f6r r\u00e4tteg\u00e5ngsmissbruk (2 kap 2 \u00a7 SkL, jfr BrB 15:2)\" `\n    \"100 000 kr \u2013 Skadest\u00e5nd f\u00f6r r\u00e4tteg\u00e5ngsmissbruk (2 kap 2 \u00a7 SkL, jfr BrB 15:2)\"\n\n# 4) Economic-damage bullet: reword + bump 25 000 kr -> 75 000 kr.\nReplace-All `\n    \"25 000 kr \u2013 Ekonomisk skada (\u00f6vers\u00e4ttningskostnader, f\u00f6rlorad arbetstid m.m.)\" `\n    \"75 000 kr \u2013 Ekonomisk skada (\u00f6vers\u00e4ttningskostnader, f\u00f6rlorad inkomst, resekostnader m.m.)\"\n\n# 5) Insert a brand-new \"Sveda och v\u00e4rk\" bullet right after the (now\n#    reworded) economic-damage bullet, before the closing sentence.\n$idx = Find-ParagraphIndex \"75 000 kr \u2013 Ekonomisk skada (\u00f6vers\u00e4ttningskostnader, f\u00f6rlorad inkomst, resekostnader m.m.)\"\nif ($idx -gt 0) {\n    $p = $d.Paragraphs.Item($idx)\n    $p.Range.InsertParagraphAfter()\n    $newPara = $d.Paragraphs.Item($idx + 1)\n    $newPara.Range.Text = \"40 000 kr \u2013 Sveda och v\u00e4rk / psykiskt lidande (5 kap 1 \u00a7 SkL)\"\n    $newPara.Style = \"List Bullet\"\n}\n\n# 6) Cross-reference filename in the closing sentence.\nReplace-All \"08_GENKAROMAL_150000kr.docx\" \"08_GENKAROMAL_400000kr.docx\"\n"}
